$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1025.4
$ws.Range("I12").Value = 825.6667
$ws.Range("J12").Value = 1325
$ws.Range("K12").Value = 825.6667
$ws.Range("L12").Value = 1325
$ws.Range("M12").Value = -655.6667
$ws.Range("N12").Value = -1665
$ws.Range("H28").Value = 58220.723
$ws.Range("J28").Value = 87047.5
$ws.Range("L28").Value = 87047.5
$ws.Range("N28").Value = -88017.5
$ws.Range("H29").Value = 3499.4
$ws.Range("J29").Value = 7999.5
$ws.Range("L29").Value = 23998.5
$ws.Range("N29").Value = -24560.5
$ws.Range("H43").Value = 5846.385
$ws.Range("I43").Value = 6401.909
$ws.Range("J43").Value = 5127.4707
$ws.Range("K43").Value = 6401.909
$ws.Range("L43").Value = 5127.4707
$ws.Range("M43").Value = -6332.909
$ws.Range("N43").Value = -5265.4707
$ws.Range("H70").Value = 1177.1471
$ws.Range("J70").Value = 1067.4333
$ws.Range("L70").Value = 3202.2999
$ws.Range("N70").Value = -3742.2999
$ws.Range("H73").Value = 1177.1471
$ws.Range("J73").Value = 1067.4333
$ws.Range("L73").Value = 3202.2999
$ws.Range("N73").Value = -5074.2999
$ws.Range("H81").Value = 199416.67
$ws.Range("J81").Value = 199416.67
$ws.Range("L81").Value = 199416.67
$ws.Range("N81").Value = -201412.67
$ws.Range("H84").Value = 199416.67
$ws.Range("J84").Value = 199416.67
$ws.Range("L84").Value = 598250.01
$ws.Range("N84").Value = -608234.01
$ws.Range("H98").Value = 594.94116
$ws.Range("J98").Value = 3200
$ws.Range("L98").Value = 3200
$ws.Range("N98").Value = -6196
$ws.Range("H106").Value = 27503664
$ws.Range("I106").Value = 31432094
$ws.Range("K106").Value = 31432094
$ws.Range("M106").Value = -31431463
$ws.Range("H107").Value = 612.13336
$ws.Range("I107").Value = 553.6923
$ws.Range("K107").Value = 553.6923
$ws.Range("M107").Value = 1366.3077
$ws.Range("H111").Value = 50194.184
$ws.Range("I111").Value = 989.3
$ws.Range("K111").Value = 2967.9
$ws.Range("M111").Value = 99.10000000000036
$ws.Range("H122").Value = 594.94116
$ws.Range("J122").Value = 3200
$ws.Range("L122").Value = 9600
$ws.Range("N122").Value = -14500
$ws.Range("H132").Value = 2955.04
$ws.Range("I132").Value = 2578.2083
$ws.Range("J132").Value = 11999
$ws.Range("K132").Value = 7734.624899999999
$ws.Range("L132").Value = 35997
$ws.Range("M132").Value = -5204.624899999999
$ws.Range("N132").Value = -41057
$ws.Range("H135").Value = 1935.9231
$ws.Range("J135").Value = 5297.5
$ws.Range("L135").Value = 47677.5
$ws.Range("N135").Value = -52747.5
$ws.Range("H137").Value = 3281.6333
$ws.Range("I137").Value = 2603.3125
$ws.Range("J137").Value = 4056.8572
$ws.Range("K137").Value = 7809.9375
$ws.Range("L137").Value = 12170.5716
$ws.Range("M137").Value = -5259.9375
$ws.Range("N137").Value = -17270.5716
$ws.Range("H138").Value = 5767.6
$ws.Range("I138").Value = 4520.727
$ws.Range("J138").Value = 7877.6924
$ws.Range("K138").Value = 13562.181
$ws.Range("L138").Value = 23633.0772
$ws.Range("M138").Value = -8422.181
$ws.Range("N138").Value = -33913.0772

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 894.9167
$ws.Range("I2").Value = 894.9167
$ws.Range("K2").Value = 894.9167
$ws.Range("M2").Value = -781.9167
$ws.Range("H32").Value = 5587.705
$ws.Range("I32").Value = 1872.3572
$ws.Range("K32").Value = 1872.3572
$ws.Range("M32").Value = -1585.3572
$ws.Range("H45").Value = 1877.0834
$ws.Range("I45").Value = 1152.7
$ws.Range("J45").Value = 5499
$ws.Range("K45").Value = 1152.7
$ws.Range("L45").Value = 5499
$ws.Range("M45").Value = -775.7
$ws.Range("N45").Value = -6253
$ws.Range("H63").Value = 2750.5386
$ws.Range("I63").Value = 1460.25
$ws.Range("J63").Value = 3856.5
$ws.Range("K63").Value = 1460.25
$ws.Range("L63").Value = 3856.5
$ws.Range("M63").Value = -774.25
$ws.Range("N63").Value = -5228.5
$ws.Range("H66").Value = 2750.5386
$ws.Range("I66").Value = 1460.25
$ws.Range("J66").Value = 3856.5
$ws.Range("K66").Value = 7301.25
$ws.Range("L66").Value = 19282.5
$ws.Range("M66").Value = -3869.25
$ws.Range("N66").Value = -26146.5
$ws.Range("H74").Value = 5206.421
$ws.Range("I74").Value = 1533
$ws.Range("J74").Value = 13165.5
$ws.Range("K74").Value = 1533
$ws.Range("L74").Value = 13165.5
$ws.Range("M74").Value = -659
$ws.Range("N74").Value = -14913.5
$ws.Range("H77").Value = 5206.421
$ws.Range("I77").Value = 1533
$ws.Range("J77").Value = 13165.5
$ws.Range("K77").Value = 7665
$ws.Range("L77").Value = 65827.5
$ws.Range("M77").Value = -3297
$ws.Range("N77").Value = -74563.5
$ws.Range("H88").Value = 2365
$ws.Range("I88").Value = 1550
$ws.Range("J88").Value = 2830.7144
$ws.Range("K88").Value = 1550
$ws.Range("L88").Value = 2830.7144
$ws.Range("M88").Value = -1144
$ws.Range("N88").Value = -3642.7144
$ws.Range("H91").Value = 2365
$ws.Range("I91").Value = 1550
$ws.Range("J91").Value = 2830.7144
$ws.Range("K91").Value = 1550
$ws.Range("L91").Value = 2830.7144
$ws.Range("M91").Value = -146
$ws.Range("N91").Value = -5638.7144
$ws.Range("H97").Value = 518.53845
$ws.Range("J97").Value = 1855.5
$ws.Range("L97").Value = 1855.5
$ws.Range("N97").Value = -2847.5
$ws.Range("H116").Value = 894.9167
$ws.Range("I116").Value = 894.9167
$ws.Range("K116").Value = 894.9167
$ws.Range("M116").Value = 1399.0833
$ws.Range("H132").Value = 52635492
$ws.Range("I132").Value = 1581.8182
$ws.Range("K132").Value = 4745.4546
$ws.Range("M132").Value = -2215.4546

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 894.9167
$ws.Range("I3").Value = 894.9167
$ws.Range("K3").Value = 894.9167
$ws.Range("M3").Value = -780.9167
$ws.Range("H94").Value = 719.7692
$ws.Range("I94").Value = 747.9524
$ws.Range("J94").Value = 601.4
$ws.Range("K94").Value = 747.9524
$ws.Range("L94").Value = 601.4
$ws.Range("M94").Value = -296.9524
$ws.Range("N94").Value = -1503.4
$ws.Range("H99").Value = 1561.8422
$ws.Range("I99").Value = 1410.3529
$ws.Range("K99").Value = 1410.3529
$ws.Range("M99").Value = 87.64709999999991
$ws.Range("H105").Value = 4082
$ws.Range("I105").Value = 6500
$ws.Range("J105").Value = 3276
$ws.Range("K105").Value = 6500
$ws.Range("L105").Value = 3276
$ws.Range("M105").Value = -4753
$ws.Range("N105").Value = -6770
$ws.Range("H134").Value = 5963.875
$ws.Range("I134").Value = 1952
$ws.Range("K134").Value = 5856
$ws.Range("M134").Value = -3321

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 384.1111
$ws.Range("I7").Value = 399.42856
$ws.Range("J7").Value = 374.36365
$ws.Range("K7").Value = 399.42856
$ws.Range("L7").Value = 374.36365
$ws.Range("M7").Value = -286.42856
$ws.Range("N7").Value = -600.36365
$ws.Range("H15").Value = 6783.3335
$ws.Range("I15").Value = 412.5
$ws.Range("J15").Value = 11880
$ws.Range("K15").Value = 412.5
$ws.Range("L15").Value = 11880
$ws.Range("M15").Value = -242.5
$ws.Range("N15").Value = -12220
$ws.Range("H16").Value = 1491.0625
$ws.Range("I16").Value = 1142.6666
$ws.Range("K16").Value = 1142.6666
$ws.Range("M16").Value = -855.6666
$ws.Range("H31").Value = 7275.0835
$ws.Range("I31").Value = 4713.5454
$ws.Range("J31").Value = 8402.16
$ws.Range("K31").Value = 4713.5454
$ws.Range("L31").Value = 8402.16
$ws.Range("M31").Value = -4418.5454
$ws.Range("N31").Value = -8992.16
$ws.Range("H34").Value = 7275.0835
$ws.Range("I34").Value = 4713.5454
$ws.Range("J34").Value = 8402.16
$ws.Range("K34").Value = 4713.5454
$ws.Range("L34").Value = 8402.16
$ws.Range("M34").Value = -4511.5454
$ws.Range("N34").Value = -8806.16
$ws.Range("H58").Value = 4898.6
$ws.Range("I58").Value = 4083.7058
$ws.Range("K58").Value = 4083.7058
$ws.Range("M58").Value = -3880.7058
$ws.Range("H107").Value = 1477.6364
$ws.Range("I107").Value = 1336.9375
$ws.Range("J107").Value = 1852.8334
$ws.Range("K107").Value = 1336.9375
$ws.Range("L107").Value = 1852.8334
$ws.Range("M107").Value = 583.0625
$ws.Range("N107").Value = -5692.8334
$ws.Range("H113").Value = 1491.0625
$ws.Range("I113").Value = 1142.6666
$ws.Range("K113").Value = 1142.6666
$ws.Range("M113").Value = 1027.3334
$ws.Range("H122").Value = 1736.8948
$ws.Range("I122").Value = 1694.1765
$ws.Range("K122").Value = 5082.529500000001
$ws.Range("M122").Value = -2632.529500000001
$ws.Range("H132").Value = 5445.1665
$ws.Range("I132").Value = 2928.5715
$ws.Range("J132").Value = 14253.25
$ws.Range("K132").Value = 8785.7145
$ws.Range("L132").Value = 42759.75
$ws.Range("M132").Value = -6255.7145
$ws.Range("N132").Value = -47819.75
$ws.Range("H134").Value = 6699.364
$ws.Range("I134").Value = 4968.3125
$ws.Range("J134").Value = 11315.5
$ws.Range("K134").Value = 14904.9375
$ws.Range("L134").Value = 33946.5
$ws.Range("M134").Value = -12369.9375
$ws.Range("N134").Value = -39016.5
$ws.Range("H136").Value = 4898.6
$ws.Range("I136").Value = 4083.7058
$ws.Range("K136").Value = 12251.1174
$ws.Range("M136").Value = -9701.117400000001
$ws.Range("H138").Value = 64789.89
$ws.Range("I138").Value = 20000
$ws.Range("J138").Value = 87184.836
$ws.Range("K138").Value = 20000
$ws.Range("L138").Value = 87184.836
$ws.Range("M138").Value = -14860
$ws.Range("N138").Value = -97464.836

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 3992.5
$ws.Range("J32").Value = 3992.5
$ws.Range("L32").Value = 11977.5
$ws.Range("N32").Value = -12543.5
$ws.Range("H62").Value = 5579.6665
$ws.Range("I62").Value = 3141.25
$ws.Range("K62").Value = 9423.75
$ws.Range("M62").Value = -8737.75
$ws.Range("H65").Value = 5579.6665
$ws.Range("I65").Value = 3141.25
$ws.Range("K65").Value = 28271.25
$ws.Range("M65").Value = -24839.25
$ws.Range("H92").Value = 233.33333
$ws.Range("J92").Value = 287.5
$ws.Range("L92").Value = 862.5
$ws.Range("N92").Value = -3358.5
$ws.Range("H106").Value = 15713.143
$ws.Range("J106").Value = 15713.143
$ws.Range("L106").Value = 47139.429
$ws.Range("N106").Value = -49031.429
$ws.Range("H107").Value = 327.23077
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 327.23077
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 981.69231
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = -4821.69231
$ws.Range("H117").Value = 923.7143
$ws.Range("I117").Value = 896.8570999999999
$ws.Range("J117").Value = 950.5714
$ws.Range("K117").Value = 2690.5713
$ws.Range("L117").Value = 2851.7142
$ws.Range("M117").Value = 751.4287000000004
$ws.Range("N117").Value = -9735.7142
$ws.Range("H128").Value = 449995.34
$ws.Range("I128").Value = 449995.34
$ws.Range("K128").Value = 1349986.02
$ws.Range("M128").Value = -1345006.02

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""
$ws.Range("H62").Value = 199500
$ws.Range("J62").Value = 199500
$ws.Range("L62").Value = 199500
$ws.Range("N62").Value = -200872
$ws.Range("H65").Value = 199500
$ws.Range("J65").Value = 199500
$ws.Range("L65").Value = 598500
$ws.Range("N65").Value = -605364
$ws.Range("H70").Value = 5640.919
$ws.Range("I70").Value = 7250
$ws.Range("J70").Value = 5596.222
$ws.Range("K70").Value = 7250
$ws.Range("L70").Value = 5596.222
$ws.Range("M70").Value = -6980
$ws.Range("N70").Value = -6136.222
$ws.Range("H73").Value = 5640.919
$ws.Range("I73").Value = 7250
$ws.Range("J73").Value = 5596.222
$ws.Range("K73").Value = 7250
$ws.Range("L73").Value = 5596.222
$ws.Range("M73").Value = -6314
$ws.Range("N73").Value = -7468.222
$ws.Range("H97").Value = 494.8
$ws.Range("I97").Value = 423.125
$ws.Range("J97").Value = 576.7143
$ws.Range("K97").Value = 423.125
$ws.Range("L97").Value = 576.7143
$ws.Range("M97").Value = 72.875
$ws.Range("N97").Value = -1568.7143
$ws.Range("H102").Value = 2738.1428
$ws.Range("I102").Value = 2189.3125
$ws.Range("J102").Value = 8592.333000000001
$ws.Range("K102").Value = 2189.3125
$ws.Range("L102").Value = 8592.333000000001
$ws.Range("M102").Value = -567.3125
$ws.Range("N102").Value = -11836.333
$ws.Range("H107").Value = 735.6
$ws.Range("J107").Value = 893.125
$ws.Range("L107").Value = 893.125
$ws.Range("N107").Value = -4733.125
$ws.Range("H110").Value = 199500
$ws.Range("J110").Value = 199500
$ws.Range("L110").Value = 199500
$ws.Range("N110").Value = -207680
$ws.Range("H113").Value = 2927.375
$ws.Range("I113").Value = 2102.3809
$ws.Range("J113").Value = 4502.364
$ws.Range("K113").Value = 2102.3809
$ws.Range("L113").Value = 4502.364
$ws.Range("M113").Value = 67.61909999999989
$ws.Range("N113").Value = -8842.364
$ws.Range("H122").Value = 7155.55
$ws.Range("I122").Value = 6101.522
$ws.Range("J122").Value = 8581.588
$ws.Range("K122").Value = 18304.566
$ws.Range("L122").Value = 25744.764
$ws.Range("M122").Value = -15854.566
$ws.Range("N122").Value = -30644.764
$ws.Range("H126").Value = 3944.1875
$ws.Range("I126").Value = 1646.091
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 4938.272999999999
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -2468.272999999999
$ws.Range("N126").Value = -31940
$ws.Range("H132").Value = 7374.75
$ws.Range("I132").Value = 2999.5
$ws.Range("K132").Value = 8998.5
$ws.Range("M132").Value = -6468.5
$ws.Range("H135").Value = 79695
$ws.Range("J135").Value = 79695
$ws.Range("L135").Value = 79695
$ws.Range("N135").Value = -89835

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = ""
$ws.Range("H22").Value = 2127.4285
$ws.Range("I22").Value = 864
$ws.Range("J22").Value = 2472
$ws.Range("K22").Value = 864
$ws.Range("L22").Value = 2472
$ws.Range("M22").Value = -569
$ws.Range("N22").Value = -3062
$ws.Range("H27").Value = 2127.4285
$ws.Range("I27").Value = 864
$ws.Range("J27").Value = 2472
$ws.Range("K27").Value = 864
$ws.Range("L27").Value = 2472
$ws.Range("M27").Value = -757
$ws.Range("N27").Value = -2686
$ws.Range("H45").Value = 14328.167
$ws.Range("I45").Value = 4000
$ws.Range("J45").Value = 19492.25
$ws.Range("K45").Value = 4000
$ws.Range("L45").Value = 19492.25
$ws.Range("M45").Value = -3593
$ws.Range("N45").Value = -20306.25
$ws.Range("H46").Value = 3762.3333
$ws.Range("I46").Value = 720
$ws.Range("J46").Value = 4522.9165
$ws.Range("K46").Value = 720
$ws.Range("L46").Value = 4522.9165
$ws.Range("M46").Value = -532
$ws.Range("N46").Value = -4898.9165
$ws.Range("H55").Value = 1348.4166
$ws.Range("I55").Value = 287.8
$ws.Range("K55").Value = 287.8
$ws.Range("M55").Value = -114.8
$ws.Range("H122").Value = 4518.579
$ws.Range("I122").Value = 3358.6924
$ws.Range("J122").Value = 7031.6665
$ws.Range("K122").Value = 10076.0772
$ws.Range("L122").Value = 21094.9995
$ws.Range("M122").Value = -7626.0772
$ws.Range("N122").Value = -25994.9995
$ws.Range("H132").Value = 10809.8
$ws.Range("I132").Value = 9489.4
$ws.Range("J132").Value = 12130.2
$ws.Range("K132").Value = 28468.2
$ws.Range("L132").Value = 36390.60000000001
$ws.Range("M132").Value = -25938.2
$ws.Range("N132").Value = -41450.60000000001
$ws.Range("H136").Value = 40006940
$ws.Range("I136").Value = 6612.4375
$ws.Range("J136").Value = 111118630
$ws.Range("K136").Value = 19837.3125
$ws.Range("L136").Value = 333355890
$ws.Range("M136").Value = -17287.3125
$ws.Range("N136").Value = -333360990

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 11450
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = ""
$ws.Range("H14").Value = 17884
$ws.Range("J14").Value = 5729.7144
$ws.Range("L14").Value = 5729.7144
$ws.Range("N14").Value = -6065.7144
$ws.Range("H38").Value = 45000
$ws.Range("J38").Value = 40000
$ws.Range("L38").Value = 40000
$ws.Range("N38").Value = -40946
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = ""
$ws.Range("N53").Value = ""
$ws.Range("H55").Value = 66666.664
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = ""
$ws.Range("H62").Value = 5003.4707
$ws.Range("I62").Value = 5663.222
$ws.Range("J62").Value = 4261.25
$ws.Range("K62").Value = 5663.222
$ws.Range("L62").Value = 4261.25
$ws.Range("M62").Value = -5039.222
$ws.Range("N62").Value = -5509.25
$ws.Range("H65").Value = 5003.4707
$ws.Range("I65").Value = 5663.222
$ws.Range("J65").Value = 4261.25
$ws.Range("K65").Value = 28316.11
$ws.Range("L65").Value = 21306.25
$ws.Range("M65").Value = -25196.11
$ws.Range("N65").Value = -27546.25
$ws.Range("H81").Value = 2333.1
$ws.Range("I81").Value = 2280.2856
$ws.Range("K81").Value = 4560.5712
$ws.Range("M81").Value = -3499.5712
$ws.Range("H84").Value = 2333.1
$ws.Range("I84").Value = 2280.2856
$ws.Range("K84").Value = 22802.856
$ws.Range("M84").Value = -17498.856
$ws.Range("H96").Value = 18022.615
$ws.Range("I96").Value = 1110
$ws.Range("K96").Value = 1110
$ws.Range("M96").Value = 263
$ws.Range("H122").Value = 2628.85
$ws.Range("I122").Value = 2567.8235
$ws.Range("J122").Value = 2974.6667
$ws.Range("K122").Value = 7703.470499999999
$ws.Range("L122").Value = 8924.000100000001
$ws.Range("M122").Value = -5253.470499999999
$ws.Range("N122").Value = -13824.0001
$ws.Range("H126").Value = 2371.625
$ws.Range("I126").Value = 2118.923
$ws.Range("J126").Value = 3466.6667
$ws.Range("K126").Value = 6356.768999999999
$ws.Range("L126").Value = 10400.0001
$ws.Range("M126").Value = -3886.768999999999
$ws.Range("N126").Value = -15340.0001
$ws.Range("H132").Value = 2465.3572
$ws.Range("I132").Value = 1895.5
$ws.Range("J132").Value = 3890
$ws.Range("K132").Value = 5686.5
$ws.Range("L132").Value = 11670
$ws.Range("M132").Value = -3156.5
$ws.Range("N132").Value = -16730
